$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "id" values in column B (rows 2-6)
$ws.Range("B2").Value = 101
$ws.Range("B3").Value = 103
$ws.Range("B4").Value = 104
$ws.Range("B5").Value = 105
$ws.Range("B6").Value = 106

# Extend column A with the next sequential ids for rows 5 and 6
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Move the active selection to A7 (single cell)
$ws.Range("A7").Select()
